$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 3557
$ws.Range("I62").Value = 1141.579
$ws.Range("K62").Value = 1141.579
$ws.Range("M62").Value = -517.579

# Row 65
$ws.Range("H65").Value = 3557
$ws.Range("I65").Value = 1141.579
$ws.Range("K65").Value = 5707.895
$ws.Range("M65").Value = -2587.895

# Row 137
$ws.Range("H137").Value = 47620572
$ws.Range("I137").Value = 58824652
$ws.Range("K137").Value = 176473956
$ws.Range("M137").Value = -176471406

# Row 138
$ws.Range("H138").Value = 3900618.5
$ws.Range("I138").Value = 1089812.5
$ws.Range("J138").Value = 5467297.5
$ws.Range("K138").Value = 3269437.5
$ws.Range("L138").Value = 16401892.5
$ws.Range("M138").Value = -3264297.5
$ws.Range("N138").Value = -16412172.5

$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 10726.917
$ws.Range("I74").Value = 1613.2
$ws.Range("J74").Value = 25916.445
$ws.Range("K74").Value = 1613.2
$ws.Range("L74").Value = 25916.445
$ws.Range("M74").Value = -739.2
$ws.Range("N74").Value = -27664.445

# Row 77
$ws.Range("H77").Value = 10726.917
$ws.Range("I77").Value = 1613.2
$ws.Range("J77").Value = 25916.445
$ws.Range("K77").Value = 8066
$ws.Range("L77").Value = 129582.225
$ws.Range("M77").Value = -3698
$ws.Range("N77").Value = -138318.225

# Row 138
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1110.1111
$ws.Range("I94").Value = 1149.5625
$ws.Range("J94").Value = 794.5
$ws.Range("K94").Value = 1149.5625
$ws.Range("L94").Value = 794.5
$ws.Range("M94").Value = -698.5625
$ws.Range("N94").Value = -1696.5

# Row 107
$ws.Range("H107").Value = 913.4545000000001
$ws.Range("I107").Value = 695.75
$ws.Range("J107").Value = 1494
$ws.Range("K107").Value = 695.75
$ws.Range("L107").Value = 1494
$ws.Range("M107").Value = 1224.25
$ws.Range("N107").Value = -5334

$ws = $wb.Worksheets.Item("CRP")
# Row 64
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

# Row 67
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

# Row 122
$ws.Range("H122").Value = 2181.7693
$ws.Range("I122").Value = 1196.3334
$ws.Range("J122").Value = 3026.4285
$ws.Range("K122").Value = 3589.0002
$ws.Range("L122").Value = 9079.2855
$ws.Range("M122").Value = -1139.0002
$ws.Range("N122").Value = -13979.2855

$ws = $wb.Worksheets.Item("CUL")
# Row 82
$ws.Range("H82").Value = 3233.1667
$ws.Range("J82").Value = 3799.8
$ws.Range("L82").Value = 11399.4
$ws.Range("N82").Value = -12211.4

# Row 85
$ws.Range("H85").Value = 3233.1667
$ws.Range("J85").Value = 3799.8
$ws.Range("L85").Value = 11399.4
$ws.Range("N85").Value = -14207.4

# Row 86
$ws.Range("H86").Value = 452.72726
$ws.Range("I86").Value = 340
$ws.Range("J86").Value = 650
$ws.Range("K86").Value = 1020
$ws.Range("L86").Value = 1950
$ws.Range("M86").Value = 166
$ws.Range("N86").Value = -4322

# Row 87
$ws.Range("H87").Value = 14922.2
$ws.Range("I87").Value = 5905.6
$ws.Range("J87").Value = 17927.732
$ws.Range("K87").Value = 17716.8
$ws.Range("L87").Value = 53783.196
$ws.Range("M87").Value = -16468.8
$ws.Range("N87").Value = -56279.196

# Row 88
$ws.Range("H88").Value = 7536.615
$ws.Range("J88").Value = 7536.615
$ws.Range("L88").Value = 22609.845
$ws.Range("N88").Value = -23465.845

# Row 89
$ws.Range("H89").Value = 452.72726
$ws.Range("I89").Value = 340
$ws.Range("J89").Value = 650
$ws.Range("K89").Value = 3060
$ws.Range("L89").Value = 5850
$ws.Range("M89").Value = 2868
$ws.Range("N89").Value = -17706

# Row 90
$ws.Range("H90").Value = 14922.2
$ws.Range("I90").Value = 5905.6
$ws.Range("J90").Value = 17927.732
$ws.Range("K90").Value = 53150.4
$ws.Range("L90").Value = 161349.588
$ws.Range("M90").Value = -46910.4
$ws.Range("N90").Value = -173829.588

# Row 91
$ws.Range("H91").Value = 7536.615
$ws.Range("J91").Value = 7536.615
$ws.Range("L91").Value = 22609.845
$ws.Range("N91").Value = -25573.845

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2775
$ws.Range("I80").Value = 2700
$ws.Range("K80").Value = 2700
$ws.Range("M80").Value = -1702

# Row 83
$ws.Range("H83").Value = 2775
$ws.Range("I83").Value = 2700
$ws.Range("K83").Value = 13500
$ws.Range("M83").Value = -8508

# Row 126
$ws.Range("H126").Value = 2322.2258
$ws.Range("I126").Value = 1708.4
$ws.Range("J126").Value = 2614.524
$ws.Range("K126").Value = 5125.200000000001
$ws.Range("L126").Value = 7843.572
$ws.Range("M126").Value = -2655.200000000001
$ws.Range("N126").Value = -12783.572

# Row 139
$ws.Range("H139").Value = 40325.668
$ws.Range("J139").Value = 40325.668
$ws.Range("L139").Value = 40325.668
$ws.Range("N139").Value = -50605.668

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3143.238
$ws.Range("I7").Value = 1725.75
$ws.Range("J7").Value = 3476.7646
$ws.Range("K7").Value = 1725.75
$ws.Range("L7").Value = 3476.7646
$ws.Range("M7").Value = -1613.75
$ws.Range("N7").Value = -3700.7646

# Row 40
$ws.Range("H40").Value = 2248.9167
$ws.Range("I40").Value = 1147.1177
$ws.Range("K40").Value = 1147.1177
$ws.Range("M40").Value = -1011.1177

# Row 122
$ws.Range("H122").Value = 3094.0881
$ws.Range("I122").Value = 1883
$ws.Range("J122").Value = 3754.682
$ws.Range("K122").Value = 5649
$ws.Range("L122").Value = 11264.046
$ws.Range("M122").Value = -3199
$ws.Range("N122").Value = -16164.046

# Row 126
$ws.Range("H126").Value = 3143.238
$ws.Range("I126").Value = 1725.75
$ws.Range("J126").Value = 3476.7646
$ws.Range("K126").Value = 5177.25
$ws.Range("L126").Value = 10430.2938
$ws.Range("M126").Value = -2707.25
$ws.Range("N126").Value = -15370.2938

$ws = $wb.Worksheets.Item("WVR")
# Row 63
$ws.Range("H63").Value = 28440
$ws.Range("J63").Value = 34300
$ws.Range("L63").Value = 34300
$ws.Range("N63").Value = -35548

# Row 66
$ws.Range("H66").Value = 28440
$ws.Range("J66").Value = 34300
$ws.Range("L66").Value = 102900
$ws.Range("N66").Value = -109140

# Row 132
$ws.Range("H132").Value = 15629032
$ws.Range("I132").Value = 21743356
$ws.Range("J132").Value = 3538.889
$ws.Range("K132").Value = 65230068
$ws.Range("L132").Value = 10616.667
$ws.Range("M132").Value = -65227538
$ws.Range("N132").Value = -15676.667
